# Add a new slide ("Title and Content" layout, the 2nd layout on the
# slide master) as the 2nd slide of the deck, give the title and body
# placeholders their text, and apply a "Fly In" entrance animation
# (effect id 2, default subtype 4 = From Bottom) to the title
# placeholder - matching the authored p:timing block that targets
# shape id 2 (the title placeholder).

$p = $ppt.ActivePresentation

# ppLayoutText (2) == the 2nd slide-layout registered on the master,
# i.e. ppt/slideLayouts/slideLayout2.xml ("Title and Content"), which
# is what slide2.xml in the target deck uses (title + idx="1" body).
$s = $p.Slides.Add(2, 2)

$s.Shapes.Item(1).TextFrame.TextRange.Text = "QWE4R5TY"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "123456"

# Entrance animation on the title placeholder: Fly In (effect id 2),
# default subtype (4 - From Bottom), triggered on click.
$seq = $s.TimeLine.MainSequence
$effect = $seq.AddEffect($s.Shapes.Item(1), 2)
